# Added tutorial level six, don't fall off the platforms.
#
# This applies:
#  - Fix typo "If play hits switch ..." -> "If player hits switch ..."
#    (C31 / C32) - the shared-string table naturally re-dedupes/renumbers.
#  - Mark three tasks "Done" with start/end dates: row 7 (Monitor ball
#    height), row 8 (Check if ball falls below a specific height) and
#    row 60 (Level 6 - can fall off level and die).
#  - Update the current selection (no more frozen topLeftCell scroll,
#    selection moved to G59).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix "play" -> "player" typo in the two switch-trigger task rows ---
$ws.Range("C31").Value = "If player hits switch animation starts/stops"
$ws.Range("C32").Value = "If player hits switch light on/off"

# --- Row 7: Monitor ball height (y coord) -> started/completed/Done ---
$ws.Range("E9").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E7").Value = 41320

$ws.Range("F9").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("F7").Value = 41320

$ws.Range("G9").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("G7").Value = "Done"

# --- Row 8: Check if ball falls below a specific height -> Done ---
$ws.Range("E9").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("E8").Value = 41320

$ws.Range("F9").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F8").Value = 41320

$ws.Range("G9").Copy()
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("G8").Value = "Done"

# --- Row 60: Level 6 - can fall off level and die -> Done ---
$ws.Range("E4").Copy()
$ws.Range("E60").PasteSpecial(-4122)
$ws.Range("E60").Value = 41320

$ws.Range("F4").Copy()
$ws.Range("F60").PasteSpecial(-4122)
$ws.Range("F60").Value = 41320

$ws.Range("G4").Copy()
$ws.Range("G60").PasteSpecial(-4122)
$ws.Range("G60").Value = "Done"

# --- Update selection / scroll position ---
$ws.Range("G59").Select()
